$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Material extra"
$ws.Range("C2").Value = "Caixas de diálogo no Shiny"
$ws.Range("D2").Value = "https://shiny.rstudio.com/articles/modal-dialogs.html"

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Material extra"
$ws.Range("C3").Value = "Pacote fresh"
$ws.Range("D3").Value = "https://github.com/dreamRs/fresh"

# Row 4 (url typed before desc)
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Material extra"
$ws.Range("D4").Value = "https://www.youtube.com/watch?v=s9GKim52E4k"
$ws.Range("C4").Value = "Vídeo: Shinydashboard, pacote fresh e CSS"

# Row 5 (aula typed before tema/desc before url, tema entered last)
$ws.Range("A5").Value = 7
$ws.Range("C5").Value = "Live sobre módulos"
$ws.Range("D5").Value = "https://www.youtube.com/watch?v=xp5aMvwqEMY&ab_channel=Curso-R"
$ws.Range("B5").Value = "Material Extra"
